$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -7.476000000000001
$ws.Range("A9").Value = -20.912
$ws.Range("D11").Value = -8.316999999999998
$ws.Range("A18").Value = -21.825
$ws.Range("A20").Value = -21.738
$ws.Range("E21").Value = 13.186
